# Update sample data definitions with better grocery figure
# (and reorder the "Income:Bonus" row up into the paycheck-related block)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Definition")

# 1) Groceries (Food:At Home / Ralph's) - lower the weekly spend figure
#    -350*52 (-18200/yr) -> -85*52 (-4420/yr)
$ws.Range("C15").Formula = "=-85*52"

# 2) Move the "Income:Bonus" / "Megacorp, Inc." row (currently row 45) so it
#    sits just above "Income:Salary" (currently row 33), pushing the
#    paycheck-withholding rows currently at 33-44 down by one to 34-45.
#
#    We do this with plain cell copies (rather than Rows.Insert/Delete) so
#    the worksheet's defined Table (A1:G46) keeps its original extent.

$bonusCategory = $ws.Range("A45").Value2
$bonusPayee    = $ws.Range("B45").Value2
$bonusAmount   = $ws.Range("C45").Value2
$bonusScheme   = $ws.Range("D45").Value2
$bonusJitterA  = $ws.Range("E45").Value2
$bonusJitterD  = $ws.Range("F45").Value2

# Shift rows 33-44 down into 34-45, working from the bottom up so a row
# isn't overwritten before it has been read.
for ($r = 44; $r -ge 33; $r--) {
    $dest = $r + 1

    $ws.Range("A$dest").Value = $ws.Range("A$r").Value2

    $bVal = $ws.Range("B$r").Value2
    if ($bVal -eq $null) {
        $ws.Range("B$dest").ClearContents() | Out-Null
    } else {
        $ws.Range("B$dest").Value = $bVal
    }

    $ws.Range("C$dest").Formula = $ws.Range("C$r").Formula
    $ws.Range("D$dest").Value = $ws.Range("D$r").Value2
    $ws.Range("E$dest").Value = $ws.Range("E$r").Value2
    $ws.Range("F$dest").Value = $ws.Range("F$r").Value2

    $gVal = $ws.Range("G$r").Value2
    if ($gVal -eq $null) {
        $ws.Range("G$dest").ClearContents() | Out-Null
    } else {
        $ws.Range("G$dest").Value = $gVal
    }
}

# Write the "Income:Bonus" row data into the now-vacated row 33 (it has no
# Group / column G value).
$ws.Range("A33").Value = $bonusCategory
$ws.Range("B33").Value = $bonusPayee
$ws.Range("C33").Value = $bonusAmount
$ws.Range("D33").Value = $bonusScheme
$ws.Range("E33").Value = $bonusJitterA
$ws.Range("F33").Value = $bonusJitterD
$ws.Range("G33").ClearContents() | Out-Null

# 3) Update the current on-screen selection to reflect where editing left off
$ws.Range("A36").Select() | Out-Null
